$wb = $excel.ActiveWorkbook

# Add the new "SearchFunction" worksheet after the existing "Filters" sheet.
$filters = $wb.Worksheets.Item("Filters")
$newSheet = $wb.Worksheets.Add($null, $filters)
$newSheet.Name = "SearchFunction"

# Populate the new sheet's data.
$newSheet.Range("A1").Value = "Assertions"
$newSheet.Range("A1").NumberFormat = "@"
$newSheet.Range("A2").Value = """DRESSES"""
$newSheet.Range("A4").Value = "PRICE DROP"
$newSheet.Range("A3").Value = "10029 Not found"

# Size column A to fit its new contents (mirrors the author's manual "best fit").
$newSheet.Columns("A:A").AutoFit()

# Update selection on the Filters sheet (no longer the active tab).
$filters.Range("A34").Select()

# Make the new sheet the active tab, with A3 selected.
$newSheet.Activate()
$newSheet.Range("A3").Select()
